$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "002"
$ws.Range("K2").Value = "001"
$ws.Range("N2").Value = "2020-06-30 00:00:00"
$ws.Range("O2").Value = 25037238.06
$ws.Range("P2").Value = 105.6060243842
$ws.Range("Q2").Value = 467585804.55
$ws.Range("R2").Value = 1972.2573935146
$ws.Range("S2").Value = 293373052.89
$ws.Range("T2").Value = 1237.4352835136
$ws.Range("U2").Value = -1841332.61
$ws.Range("V2").Value = -7.7666640404
$ws.Range("W2").Value = ""
$ws.Range("X2").Value = ""
$ws.Range("Y2").Value = 1858217.48
$ws.Range("Z2").Value = 7.8378837168
$ws.Range("AA2").Value = 498040
$ws.Range("AB2").Value = 2.1007119179
$ws.Range("AC2").Value = 23708153.21
$ws.Range("AD2").Value = ""
